# Applies the "Smaller Problems/Constraints", "What isn't visible from word
# problem", and "Possible Solutions" sections to the end of the "Socks In
# The Dark" problem, mirroring the Problem #2 sections added for the other
# word problems earlier in the document.

$d = $word.ActiveDocument

# The existing "_GoBack" bookmark currently trails the very last paragraph
# ("... pair of each color sock with the minimum amount of tries."). It
# needs to move to the end of the new "Once decision is made..." paragraph
# being appended below, so drop it here; the new block re-adds it (as raw
# markup) in its new home.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$lastPara = $d.Paragraphs.Last
$endPos = $lastPara.Range.End
$insertionPoint = $d.Range($endPos, $endPos)

$curlyQuote = [char]0x2019

$body = ''
$body += '<w:p/>'
$body += '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr>'
$body += '<w:r><w:rPr><w:b/></w:rPr><w:t>Smaller Problems/Constraints:</w:t></w:r>'
$body += '</w:p>'
$body += '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr></w:p>'
$body += '<w:p><w:r><w:t>It is dark making selecting a pair of matching socks very difficult.</w:t></w:r></w:p>'
$body += '<w:p><w:r><w:t>I am unsure of my selection as I can' + $curlyQuote + 't see what I am picking and doubt my decisions.</w:t></w:r></w:p>'
$body += '<w:p>'
$body += '<w:r><w:t xml:space="preserve">I have no real knowledge on </w:t></w:r>'
$body += '<w:r><w:t>probability, which</w:t></w:r>'
$body += '<w:r><w:t xml:space="preserve"> influences my doubts. </w:t></w:r>'
$body += '</w:p>'
$body += '<w:p>'
$body += '<w:r><w:t xml:space="preserve">Once decision is made, I run the chance of having </w:t></w:r>'
$body += '<w:proofErr w:type="spellStart"/>'
$body += '<w:r><w:t>mis</w:t></w:r>'
$body += '<w:proofErr w:type="spellEnd"/>'
$body += '<w:r><w:t>-matching socks in different colors.</w:t></w:r>'
$body += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$body += '</w:p>'
$body += '<w:p/>'
$body += '<w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr>'
$body += '<w:r><w:rPr><w:b/></w:rPr><w:t>What isn' + $curlyQuote + 't visible from word problem:</w:t></w:r>'
$body += '</w:p>'
$body += '<w:p/>'
$body += '<w:p><w:r><w:t>What is the reason for me being in the dark?</w:t></w:r></w:p>'
$body += '<w:p/>'
$body += '<w:p/>'
$body += '<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Possible Solutions:</w:t></w:r></w:p>'
$body += '<w:p/>'
$body += '<w:p/>'

$xmlFrag = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' `
    + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
    + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
    + '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' `
    + $body `
    + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($xmlFrag)

Write-Host "Problem #2 sections appended."
